$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9142543333333334
$ws.Range("H2").Value = 2.742763
$ws.Range("I2").Value = 0.1175834869881751
$ws.Range("J2").Value = 0.1175834869881751
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 153.6951381154184
$ws.Range("R2").Value = 1383.256243038766
$ws.Range("S2").Value = 0.03508905557077097
$ws.Range("T2").Value = 0.03508905557077097

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9142543333333334
$ws.Range("H3").Value = 2.742763
$ws.Range("I3").Value = 0.1175834869881751
$ws.Range("J3").Value = 0.1175834869881751
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 149.0291573186046
$ws.Range("R3").Value = 1341.262415867441
$ws.Range("S3").Value = 0.03402379832529712
$ws.Range("T3").Value = 0.03402379832529712

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9142543333333334
$ws.Range("H4").Value = 2.742763
$ws.Range("I4").Value = 0.1175834869881751
$ws.Range("J4").Value = 0.1175834869881751
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 151.7603041077967
$ws.Range("R4").Value = 1365.84273697017
$ws.Range("S4").Value = 0.03464732723215121
$ws.Range("T4").Value = 0.03464732723215121

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9142543333333334
$ws.Range("H5").Value = 2.742763
$ws.Range("I5").Value = 0.1175834869881751
$ws.Range("J5").Value = 0.1175834869881751
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 60.54807884676546
$ws.Range("R5").Value = 544.932709620889
$ws.Range("S5").Value = 0.01382330585995577
$ws.Range("T5").Value = 0.01382330585995577

$ws.Range("I6").Value = 0.6206849497708361
$ws.Range("J6").Value = 0.620684949770836
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 811.3066003118606
$ws.Range("R6").Value = 7301.759402806745
$ws.Range("S6").Value = 0.1852237014933935
$ws.Range("T6").Value = 0.1852237014933935

$ws.Range("I7").Value = 0.6206849497708361
$ws.Range("J7").Value = 0.620684949770836
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.179600555277578
$ws.Range("T7").Value = 0.1796005552775779

$ws.Range("I8").Value = 0.6206849497708361
$ws.Range("J8").Value = 0.620684949770836
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 801.09324145003
$ws.Range("R8").Value = 7209.83917305027
$ws.Range("S8").Value = 0.1828919613937302
$ws.Range("T8").Value = 0.1828919613937302

$ws.Range("I9").Value = 0.6206849497708361
$ws.Range("J9").Value = 0.620684949770836
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 319.6135974561177
$ws.Range("R9").Value = 2876.522377105059
$ws.Range("S9").Value = 0.07296873160613447
$ws.Range("T9").Value = 0.07296873160613446

$ws.Range("G10").Value = 1.986145
$ws.Range("H10").Value = 5.958435
$ws.Range("I10").Value = 0.2554407961214246
$ws.Range("J10").Value = 0.2554407961214246
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 333.8904930089633
$ws.Range("R10").Value = 3005.01443708067
$ws.Range("S10").Value = 0.07622818917632572
$ws.Range("T10").Value = 0.07622818917632573

$ws.Range("G11").Value = 1.986145
$ws.Range("H11").Value = 5.958435
$ws.Range("I11").Value = 0.2554407961214246
$ws.Range("J11").Value = 0.2554407961214246
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 323.7540199381717
$ws.Range("R11").Value = 2913.786179443545
$ws.Range("S11").Value = 0.07391400233063948
$ws.Range("T11").Value = 0.07391400233063948

$ws.Range("G12").Value = 1.986145
$ws.Range("H12").Value = 5.958435
$ws.Range("I12").Value = 0.2554407961214246
$ws.Range("J12").Value = 0.2554407961214246
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 329.6872196418499
$ws.Range("R12").Value = 2967.18497677665
$ws.Range("S12").Value = 0.07526856940847709
$ws.Range("T12").Value = 0.07526856940847709

$ws.Range("G13").Value = 1.986145
$ws.Range("H13").Value = 5.958435
$ws.Range("I13").Value = 0.2554407961214246
$ws.Range("J13").Value = 0.2554407961214246
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 131.5358972624783
$ws.Range("R13").Value = 1183.823075362305
$ws.Range("S13").Value = 0.03003003520598227
$ws.Range("T13").Value = 0.03003003520598227

$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.04891300000000001
$ws.Range("H14").Value = 0.146739
$ws.Range("I14").Value = 0.006290767119564404
$ws.Range("J14").Value = 0.006290767119564403
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 8.222755984355334
$ws.Range("R14").Value = 74.00480385919801
$ws.Range("S14").Value = 0.001877279562761843
$ws.Range("T14").Value = 0.001877279562761843

$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.04891300000000001
$ws.Range("H15").Value = 0.146739
$ws.Range("I15").Value = 0.006290767119564404
$ws.Range("J15").Value = 0.006290767119564403
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 7.973124005163668
$ws.Range("R15").Value = 71.75811604647301
$ws.Range("S15").Value = 0.001820287842024913
$ws.Range("T15").Value = 0.001820287842024912

$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.04891300000000001
$ws.Range("H16").Value = 0.146739
$ws.Range("I16").Value = 0.006290767119564404
$ws.Range("J16").Value = 0.006290767119564403
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 8.119241532890001
$ws.Range("R16").Value = 73.07317379601001
$ws.Range("S16").Value = 0.001853646906684477
$ws.Range("T16").Value = 0.001853646906684477

$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.04891300000000001
$ws.Range("H17").Value = 0.146739
$ws.Range("I17").Value = 0.006290767119564404
$ws.Range("J17").Value = 0.006290767119564403
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 3.239348256446334
$ws.Range("R17").Value = 29.154134308017
$ws.Range("S17").Value = 0.0007395528080931709
$ws.Range("T17").Value = 0.0007395528080931708
